# Apply targeted cell value updates described by the source diff.
# Odds/lay values changed across rows 2-21 of Sheet1; all other cells are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.12
$ws.Range("N2").Value = 2.62
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 1.54
$ws.Range("Q2").Value = 2.58
$ws.Range("R2").Value = 1.19
$ws.Range("S2").Value = 5.3
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 1.72
$ws.Range("X2").Value = 8.800000000000001
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 970
$ws.Range("AA2").Value = 190
$ws.Range("AB2").Value = 6.6
$ws.Range("AC2").Value = 7.6
$ws.Range("AD2").Value = 970
$ws.Range("AE2").Value = 120
$ws.Range("AF2").Value = 11
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 970
$ws.Range("AI2").Value = 140
$ws.Range("AJ2").Value = 970
$ws.Range("AK2").Value = 970
$ws.Range("AL2").Value = 60
$ws.Range("AM2").Value = 250
$ws.Range("AN2").Value = 970
$ws.Range("AO2").Value = 200
# Row 3
$ws.Range("F3").Value = 2.62
$ws.Range("G3").Value = 2.98
$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 3.2
$ws.Range("K3").Value = 3.5
$ws.Range("L3").Value = 1.48
$ws.Range("M3").Value = 1.09
$ws.Range("N3").Value = 2.9
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 1.65
$ws.Range("Q3").Value = 2.22
$ws.Range("R3").Value = 1.24
$ws.Range("S3").Value = 4.2
$ws.Range("T3").Value = 1.89
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.45
$ws.Range("W3").Value = 1.5
$ws.Range("X3").Value = 12.5
$ws.Range("Y3").Value = 10.5
$ws.Range("Z3").Value = 22
$ws.Range("AA3").Value = 60
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 8.4
$ws.Range("AD3").Value = 15.5
$ws.Range("AE3").Value = 46
$ws.Range("AF3").Value = 21
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 65
$ws.Range("AJ3").Value = 55
$ws.Range("AK3").Value = 42
$ws.Range("AL3").Value = 65
$ws.Range("AM3").Value = 150
$ws.Range("AN3").Value = 44
$ws.Range("AO3").Value = 50
# Row 4
$ws.Range("G4").Value = 3.45
$ws.Range("I4").Value = 970
$ws.Range("J4").Value = 2.38
$ws.Range("N4").Value = 1.89
$ws.Range("P4").Value = 1.89
$ws.Range("S4").Value = 2.6
$ws.Range("V4").Value = 1.01
# Row 5
$ws.Range("G5").Value = 18.5
$ws.Range("V5").Value = 5
$ws.Range("AC5").Value = 24
# Row 7
$ws.Range("F7").Value = 3.55
# Row 8
$ws.Range("I8").Value = 2.22
$ws.Range("Q8").Value = 1.59
$ws.Range("V8").Value = 1.82
# Row 9
$ws.Range("M9").Value = 1.17
$ws.Range("N9").Value = 2.18
$ws.Range("O9").Value = 1.75
$ws.Range("P9").Value = 1.38
$ws.Range("Q9").Value = 3.25
$ws.Range("T9").Value = 2.62
$ws.Range("U9").Value = 1.53
$ws.Range("W9").Value = 1.94
$ws.Range("X9").Value = 6.6
$ws.Range("Y9").Value = 14.5
$ws.Range("Z9").Value = 55
$ws.Range("AA9").Value = 280
$ws.Range("AB9").Value = 5.5
$ws.Range("AC9").Value = 9.199999999999999
$ws.Range("AD9").Value = 34
$ws.Range("AE9").Value = 190
$ws.Range("AF9").Value = 12
$ws.Range("AG9").Value = 15.5
$ws.Range("AH9").Value = 44
$ws.Range("AI9").Value = 240
$ws.Range("AJ9").Value = 32
$ws.Range("AK9").Value = 44
$ws.Range("AL9").Value = 120
$ws.Range("AM9").Value = 490
$ws.Range("AN9").Value = 42
$ws.Range("AO9").Value = 420
# Row 10
$ws.Range("L10").Value = 1.45
$ws.Range("Q10").Value = 2.1
$ws.Range("AC10").Value = 8.4
# Row 11
$ws.Range("AL11").Value = 48
# Row 12
$ws.Range("F12").Value = 1.83
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = 4.9
$ws.Range("S12").Value = 2.88
$ws.Range("T12").Value = 1.69
# Row 13
$ws.Range("G13").Value = 2.7
$ws.Range("I13").Value = 3.5
$ws.Range("V13").Value = 1.42
$ws.Range("W13").Value = 1.59
$ws.Range("AC13").Value = 9
# Row 14
$ws.Range("X14").Value = 970
$ws.Range("Y14").Value = 970
$ws.Range("AD14").Value = 970
$ws.Range("AH14").Value = 970
$ws.Range("AL14").Value = 970
# Row 15
$ws.Range("R15").Value = 1.35
# Row 16
$ws.Range("F16").Value = 9.199999999999999
$ws.Range("G16").Value = 9.6
$ws.Range("Y16").Value = 15.5
$ws.Range("AF16").Value = 95
$ws.Range("AK16").Value = 110
$ws.Range("AL16").Value = 80
# Row 17
$ws.Range("F17").Value = 3.8
$ws.Range("H17").Value = 2.06
$ws.Range("I17").Value = 2.08
$ws.Range("J17").Value = 3.9
$ws.Range("K17").Value = 3.95
$ws.Range("V17").Value = 1.92
$ws.Range("AF17").Value = 29
$ws.Range("AJ17").Value = 70
$ws.Range("AN17").Value = 29
# Row 18
$ws.Range("T18").Value = 1.91
# Row 19
$ws.Range("G19").Value = 1.51
$ws.Range("H19").Value = 6.4
$ws.Range("S19").Value = 2.04
$ws.Range("W19").Value = 2.98
$ws.Range("AO19").Value = 44
# Row 20
$ws.Range("T20").Value = 2.16
$ws.Range("Y20").Value = 85
# Row 21
$ws.Range("H21").Value = 26
$ws.Range("N21").Value = 12
$ws.Range("R21").Value = 2.44
$ws.Range("U21").Value = 1.89
$ws.Range("W21").Value = 8.800000000000001
$ws.Range("X21").Value = 990
$ws.Range("AJ21").Value = 9.6
$ws.Range("AL21").Value = 40
